# Applies the crypto price/volume refresh described in the commit diff.
# Numeric-looking "Price" values are apostrophe-prefixed so Excel stores them
# as literal text (matching the original inlineStr cells) instead of coercing
# them to floating-point numbers; ClearFormats() then drops the incidental
# quote-prefix style Excel attaches so the cell keeps its original (default)
# formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.553.92'
$ws.Range('E2').Value = '  +1.48%  '
$ws.Range('D3').Value = '1.905.02'
$ws.Range('E3').Value = '  +3.16%  '
$ws.Range('E4').Value = '  +0.55%  '
$ws.Range('D5').Value = "'247.52"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +6.11%  '
$ws.Range('D6').Value = "'0.633"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.73%  '
$ws.Range('E7').Value = '  +0.54%  '
$ws.Range('D8').Value = "'42.19"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  +3.18%  '
$ws.Range('D10').Value = "'0.0705"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +1.50%  '
$ws.Range('D11').Value = "'0.0998"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.48%  '
$ws.Range('D12').Value = '2.180.02'
$ws.Range('E12').Value = '  +3.06%  '
$ws.Range('D13').Value = "'12.41"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +9.31%  '
$ws.Range('D14').Value = '1.915.61'
$ws.Range('E14').Value = '  +3.55%  '
$ws.Range('D15').Value = "'0.691"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.14%  '
$ws.Range('D16').Value = "'4.86"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +3.60%  '
$ws.Range('D17').Value = '35.533.71'
$ws.Range('E17').Value = '  +1.35%  '
$ws.Range('D18').Value = "'72.07"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.01%  '
$ws.Range('D19').Value = '0.0₃0815'
$ws.Range('E19').Value = '  +2.73%  '
$ws.Range('D20').Value = "'244.12"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.31%  '
$ws.Range('D21').Value = "'12.50"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.15%  '
$ws.Range('D22').Value = "'4.89"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.20%  '
$ws.Range('E23').Value = '  +0.57%  '
$ws.Range('D24').Value = "'2.28"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('D25').Value = "'2.25"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +27.30%  '
$ws.Range('D26').Value = "'172.04"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.28%  '
$ws.Range('D27').Value = "'8.60"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +9.11%  '
$ws.Range('D28').Value = "'17.99"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.25%  '
$ws.Range('E29').Value = '  +1.01%  '
$ws.Range('D30').Value = "'0.972"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +27.01%  '
$ws.Range('E31').Value = '  +3.40%  '
$ws.Range('D32').Value = "'0.0567"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.04%  '
$ws.Range('E33').Value = '  +0.62%  '
$ws.Range('D34').Value = "'4.19"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +5.30%  '
$ws.Range('E35').Value = '  +5.86%  '
$ws.Range('E36').Value = '  +2.27%  '
$ws.Range('E37').Value = '  +7.30%  '
$ws.Range('D38').Value = "'1.11"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +3.02%  '
$ws.Range('E39').Value = '  +2.32%  '
$ws.Range('D40').Value = "'91.06"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.48%  '
$ws.Range('D41').Value = '1.356.13'
$ws.Range('E41').Value = '  +0.62%  '
$ws.Range('D42').Value = "'15.57"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +6.30%  '
$ws.Range('D43').Value = "'48.90"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +41.48%  '
$ws.Range('D44').Value = "'0.0593"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +11.55%  '
$ws.Range('E45').Value = '  +1.92%  '
$ws.Range('D46').Value = "'12.60"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.72%  '
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('D48').Value = "'6.67"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +5.02%  '
$ws.Range('D49').Value = "'2.76"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('D50').Value = '2.091.62'
$ws.Range('E50').Value = '  +3.11%  '
$ws.Range('D51').Value = "'0.0692"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.64%  '
